# Update cryptocurrency price/volume table (D: Price, E: Volume(1h))
# D-column values are written as literal text (matching the source feed's
# formatting, e.g. preserved trailing zeros / locale thousands separators),
# so we force NumberFormat to Text before assignment and restore the default
# "Normal" style afterwards to avoid leaving any formatting side effects.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "34.181.73"
$ws.Range("E2").Value = "  +0.91%  "
Set-TextValue "D3" "1.780.19"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.26%  "
Set-TextValue "D5" "226.06"
$ws.Range("E5").Value = "  +0.88%  "
Set-TextValue "D6" "0.546"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.27%  "
Set-TextValue "D8" "31.73"
$ws.Range("E8").Value = "  -0.37%  "
Set-TextValue "D9" "0.291"
$ws.Range("E9").Value = "  +0.70%  "
Set-TextValue "D10" "0.0690"
$ws.Range("E10").Value = "  +1.82%  "
Set-TextValue "D11" "0.0946"
$ws.Range("E11").Value = "  +1.03%  "
Set-TextValue "D12" "2.041.04"
$ws.Range("E12").Value = "  +0.22%  "
Set-TextValue "D13" "10.94"
$ws.Range("E13").Value = "  -2.40%  "
Set-TextValue "D14" "1.771.02"
$ws.Range("E14").Value = "  -0.43%  "
Set-TextValue "D15" "34.163.34"
$ws.Range("E15").Value = "  +0.85%  "
Set-TextValue "D16" "0.623"
$ws.Range("E16").Value = "  +2.18%  "
Set-TextValue "D17" "4.17"
$ws.Range("E17").Value = "  +0.82%  "
Set-TextValue "D18" "67.83"
$ws.Range("E18").Value = "  +1.66%  "
Set-TextValue "D19" "0.0₃0799"
$ws.Range("E19").Value = "  +3.49%  "
Set-TextValue "D20" "245.88"
$ws.Range("E20").Value = "  +2.94%  "
Set-TextValue "D21" "10.95"
$ws.Range("E21").Value = "  +3.43%  "
$ws.Range("E22").Value = "  +0.16%  "
Set-TextValue "D23" "4.08"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  -1.38%  "
Set-TextValue "D25" "162.12"
$ws.Range("E25").Value = "  +0.87%  "
Set-TextValue "D26" "7.18"
$ws.Range("E26").Value = "  +2.00%  "
Set-TextValue "D27" "16.27"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +2.08%  "
Set-TextValue "D32" "3.73"
$ws.Range("E32").Value = "  +4.00%  "
Set-TextValue "D33" "3.72"
$ws.Range("E33").Value = "  +5.57%  "
$ws.Range("E34").Value = "  -1.24%  "
Set-TextValue "D35" "1.439.74"
$ws.Range("E35").Value = "  +3.52%  "
Set-TextValue "D36" "0.653"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("E39").Value = "  +0.22%  "
Set-TextValue "D40" "80.10"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("E41").Value = "  -0.32%  "
Set-TextValue "D42" "0.921"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  +0.74%  "
Set-TextValue "D44" "13.50"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +0.17%  "
Set-TextValue "D46" "6.07"
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -4.15%  "
Set-TextValue "D49" "1.942.48"
$ws.Range("E49").Value = "  +0.23%  "
Set-TextValue "D50" "104.37"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  +0.24%  "

Write-Output "Updated cryptos list"
